$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new (blank) rows at these positions, in ascending order so each
# insertion's row number is relative to the sheet as it stands at that
# moment (inserting low-to-high keeps later target rows correct because
# everything below an insertion point shifts down by one row).
$insertAt = @(1, 5, 7, 14, 20)
foreach ($r in $insertAt) {
    $ws.Rows.Item($r).Insert()
    # Insert() carries column A's number-format down from the row above;
    # strip it again so the new separator row is genuinely blank.
    $ws.Rows.Item($r).ClearFormats()
}

# Row 1 becomes the new header row: Time / kind / num
$ws.Range("A1").Value = "Time"
$ws.Range("B1").Value = "kind"
$ws.Range("C1").Value = "num"
